# Update Name of Algo
# Applies updated imputed values produced by the RandomForest algorithm
# to the corresponding cells in columns A and E of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 15.05849999999999
$ws.Range("A10").Value = -22.00759999999999
$ws.Range("A12").Value = -21.48300000000001
$ws.Range("E15").Value = 15.96990000000001
$ws.Range("A18").Value = -22.2388
$ws.Range("E20").Value = 15.95489999999999
$ws.Range("E29").Value = 17.09840000000002
$ws.Range("E30").Value = 15.6125
$ws.Range("E31").Value = 16.06750000000001
$ws.Range("A37").Value = -19.7972
$ws.Range("E40").Value = 17.0018
$ws.Range("A55").Value = -22.28540000000001
$ws.Range("A68").Value = -21.50739999999999
$ws.Range("E68").Value = 16.97390000000001
$ws.Range("E76").Value = 16.27839999999999
$ws.Range("A77").Value = -20.38559999999999
$ws.Range("A78").Value = -19.87589999999998
$ws.Range("E87").Value = 16.07929999999999
$ws.Range("E88").Value = 16.30380000000001
$ws.Range("E96").Value = 16.16919999999999
$ws.Range("E98").Value = 15.4873
$ws.Range("E101").Value = 16.94970000000001
$ws.Range("E102").Value = 16.70040000000002
